$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "232.28") that Excel
# would otherwise auto-convert to a real number. Force text format before
# writing so the values round-trip as strings (matching the source data),
# then restore the original "Normal" style so no stray formatting sticks.
$textCells = @("D5","D8","D15","D18","D20","D21","D25","D29","D34","D35","D39","D42","D45","D46","D47","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.921.04"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").Value = "1.839.88"
$ws.Range("E3").Value = "  +1.85%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "232.28"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("E6").Value = "  +2.35%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "41.02"
$ws.Range("E8").Value = "  +4.59%  "

$ws.Range("E9").Value = "  +2.95%  "

$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("D12").Value = "2.107.20"

$ws.Range("E13").Value = "  +4.13%  "

$ws.Range("D14").Value = "1.839.06"
$ws.Range("E14").Value = "  +1.96%  "

$ws.Range("D15").Value = "0.669"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("E16").Value = "  +2.40%  "

$ws.Range("D17").Value = "34.937.05"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").Value = "69.81"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("D20").Value = "239.89"
$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("D21").Value = "12.13"
$ws.Range("E21").Value = "  +3.04%  "

$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("D25").Value = "172.00"
$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("E26").Value = "  +2.03%  "

$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("E28").Value = "  +3.66%  "

$ws.Range("D29").Value = "1.67"
$ws.Range("E29").Value = "  +10.66%  "

$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("E31").Value = "  +1.22%  "

$ws.Range("E32").Value = "  -0.95%  "

$ws.Range("E33").Value = "  -0.82%  "

$ws.Range("D34").Value = "1.63"
$ws.Range("E34").Value = "  +22.77%  "

$ws.Range("D35").Value = "1.94"
$ws.Range("E35").Value = "  +10.57%  "

$ws.Range("E36").Value = "  -1.19%  "

$ws.Range("E37").Value = "  +8.05%  "

$ws.Range("E38").Value = "  +10.33%  "

$ws.Range("D39").Value = "89.65"
$ws.Range("E39").Value = "  -1.50%  "

$ws.Range("E40").Value = "  +3.36%  "

$ws.Range("D41").Value = "1.337.15"
$ws.Range("E41").Value = "  +2.14%  "

$ws.Range("D42").Value = "14.59"
$ws.Range("E42").Value = "  +2.72%  "

$ws.Range("E43").Value = "  -2.16%  "

$ws.Range("E44").Value = "  +1.70%  "

$ws.Range("D45").Value = "2.74"
$ws.Range("E45").Value = "  +3.38%  "

$ws.Range("D46").Value = "0.0530"
$ws.Range("E46").Value = "  +3.86%  "

$ws.Range("D47").Value = "6.31"
$ws.Range("E47").Value = "  +2.96%  "

$ws.Range("D48").Value = "2.029.00"
$ws.Range("E48").Value = "  +1.75%  "

$ws.Range("D49").Value = "11.00"
$ws.Range("E49").Value = "  +68.12%  "

$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("D51").Value = "3.38"
$ws.Range("E51").Value = "  +15.50%  "

# Restore the default "Normal" style on the text-forced cells so the
# temporary NumberFormat="@" doesn't leave stray formatting behind.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
